$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'75.352.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.82%  "
$ws.Range("D3").Value = "'2.813.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +13.19%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "'189.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.31%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'602.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.03%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.544"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.53%  "
$ws.Range("D9").Value = "'0.199"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.20%  "
$ws.Range("D10").Value = "'2.810.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.27%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "'0.373"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.37%  "
$ws.Range("D13").Value = "'4.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("D14").Value = "'3.327.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +13.12%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000193"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.30%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'75.279.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.87%  "
$ws.Range("D17").Value = "'27.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +11.08%  "
$ws.Range("D18").Value = "'2.800.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +12.72%  "
$ws.Range("D19").Value = "'9.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +24.06%  "
$ws.Range("D20").Value = "'12.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.72%  "
$ws.Range("D21").Value = "'382.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.43%  "
$ws.Range("D22").Value = "'2.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.81%  "
$ws.Range("D23").Value = "'4.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.21%  "
$ws.Range("D24").Value = "'6.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.49%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'71.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.02%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'4.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.60%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "'2.960.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.31%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'9.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.59%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0000106"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.59%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "'529.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.73%  "
$ws.Range("D33").Value = "'1.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.43%  "
$ws.Range("D34").Value = "'7.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.35%  "
$ws.Range("D35").Value = "'1.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.59%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'20.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.95%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.121"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.03%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'161.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'182.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +28.04%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'5.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.37%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.344"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.81%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.87%  "
$ws.Range("D46").Value = "'1.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.78%  "
$ws.Range("D47").Value = "'2.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.31%  "
$ws.Range("D48").Value = "'39.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.55%  "
$ws.Range("D49").Value = "'0.0860"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.98%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'0.574"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.45%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'3.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.39%  "
